$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet to match the new source (PubMed affiliation filter)
$ws.Name = "utrc_active_allocations"

# Header row: same institution column, but last_name/first_name swapped
$ws.Range("A1").Value = "root_institution_name"
$ws.Range("B1").Value = "last_name"
$ws.Range("C1").Value = "first_name"

# Row 2 - Kelsey Beavers, UT Austin
$ws.Range("A2").Value = "The University of Texas at Austin"
$ws.Range("B2").Value = "Beavers"
$ws.Range("C2").Value = "Kelsey m"

# Row 3 - James Carson, UT Austin
$ws.Range("A3").Value = "The University of Texas at Austin"
$ws.Range("B3").Value = "Carson"
$ws.Range("C3").Value = "James"

# Row 4 - Laura Mydlarz, UT Arlington (last name/first name entered before institution)
$ws.Range("B4").Value = "Mydlarz"
$ws.Range("C4").Value = "Laura"
$ws.Range("A4").Value = "The University of Texas at Arlington"

# The old sheet had a 5th data row that no longer exists
$ws.Range("A5:C5").Delete() | Out-Null

# New data rows use the workbook's default font explicitly (Aptos Narrow 12)
$ws.Range("A2:C4").Font.Name = "Aptos Narrow"
$ws.Range("A2:C4").Font.Size = 12

# Column A is widened (auto-fit) to show the full institution name
$ws.Columns.Item(1).ColumnWidth = 26

# Mirror the saved selection/view state
$ws.Range("A5:XFD8").Select() | Out-Null
